$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("samples_retained")

# Fill in row 27 with the new "urdu" dataset entry (columns B-K)
$ws.Range("B27").Value = "discourse"
$ws.Range("C27").Value = "spon."
$ws.Range("D27").Value = 100
$ws.Range("E27").Value = 200
$ws.Range("F27").Value = 100
$ws.Range("G27").Value = "Urdu"
$ws.Range("H27").Formula = "=IF(OR(ISBLANK(D27), ISBLANK(E27),ISBLANK(F27)), """", SUM(D27:F27))"
$ws.Range("I27").Value = 29
$ws.Range("J27").Value = "angry, happy, neutral, sad"
$ws.Range("K27").Value = "expected 38 unique speakers, but only found 29; but the expected number of utterances was found"

# Update the view selection/scroll position to match the saved workbook state
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("B28").Select()
